$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# ---- Sheet: LP1912 ----
$data = New-Object 'object[,]' 2,5
$data[0,0] = 'Última actualización: 14:45:17'
$data[1,0] = 'Total filas: 292'
$ws1.Range("A2:E3").Value = $data

$data = New-Object 'object[,]' 2,5
$data[0,0] = '05:26:08'
$data[0,1] = '07:05'
$data[0,2] = '23_HERNANDEZ'
$data[0,3] = 99
$data[0,4] = 'LP1912'
$data[1,0] = '05:26:08'
$data[1,1] = '07:05'
$data[1,2] = '15_ABASTO'
$data[1,3] = 99
$data[1,4] = 'LP1912'
$ws1.Range("A40:E41").Value = $data

$data = New-Object 'object[,]' 3,5
$data[0,0] = '08:02:22'
$data[0,1] = '09:23'
$data[0,2] = '11_ETCHEVERRY'
$data[0,3] = 81
$data[0,4] = 'LP1912'
$data[1,0] = '07:50:16'
$data[1,1] = '09:23'
$data[1,2] = '17_ROMERO'
$data[1,3] = 93
$data[1,4] = 'LP1912'
$data[2,0] = '08:32:09'
$data[2,1] = '09:23'
$data[2,2] = '16_SANTA ANA'
$data[2,3] = 51
$data[2,4] = 'LP1912'
$ws1.Range("A103:E105").Value = $data

$data = New-Object 'object[,]' 2,5
$data[0,0] = '08:32:09'
$data[0,1] = '09:35'
$data[0,2] = '16_SANTA ANA'
$data[0,3] = 63
$data[0,4] = 'LP1912'
$data[1,0] = '08:48:08'
$data[1,1] = '09:35'
$data[1,2] = '23_HERNANDEZ'
$data[1,3] = 47
$data[1,4] = 'LP1912'
$ws1.Range("A112:E113").Value = $data

$data = New-Object 'object[,]' 2,5
$data[0,0] = '09:35:26'
$data[0,1] = '11:06'
$data[0,2] = '16_P MOR-167 Y 521'
$data[0,3] = 91
$data[0,4] = 'LP1912'
$data[1,0] = '10:59:49'
$data[1,1] = '11:06'
$data[1,2] = '23_HERNANDEZ'
$data[1,3] = 7
$data[1,4] = 'LP1912'
$ws1.Range("A150:E151").Value = $data

$data = New-Object 'object[,]' 3,5
$data[0,0] = '10:59:49'
$data[0,1] = '12:20'
$data[0,2] = '14_ABASTO'
$data[0,3] = 81
$data[0,4] = 'LP1912'
$data[1,0] = '10:29:57'
$data[1,1] = '12:20'
$data[1,2] = '26_HERNANDEZ'
$data[1,3] = 111
$data[1,4] = 'LP1912'
$data[2,0] = '10:29:57'
$data[2,1] = '12:20'
$data[2,2] = '215A_EL PATO'
$data[2,3] = 111
$data[2,4] = 'LP1912'
$ws1.Range("A187:E189").Value = $data

$data = New-Object 'object[,]' 2,5
$data[0,0] = '12:21:08'
$data[0,1] = '12:35'
$data[0,2] = '23_HERNANDEZ'
$data[0,3] = 14
$data[0,4] = 'LP1912'
$data[1,0] = '12:21:08'
$data[1,1] = '12:35'
$data[1,2] = '11_ETCHEVERRY'
$data[1,3] = 14
$data[1,4] = 'LP1912'
$ws1.Range("A196:E197").Value = $data

$data = New-Object 'object[,]' 2,5
$data[0,0] = '11:30:45'
$data[0,1] = '13:26'
$data[0,2] = '14_ABASTO'
$data[0,3] = 116
$data[0,4] = 'LP1912'
$data[1,0] = '11:30:45'
$data[1,1] = '13:26'
$data[1,2] = '15_ABASTO'
$data[1,3] = 116
$data[1,4] = 'LP1912'
$ws1.Range("A224:E225").Value = $data

$data = New-Object 'object[,]' 2,5
$data[0,0] = '13:33:42'
$data[0,1] = '13:34'
$data[0,2] = '16_SANTA ANA'
$data[0,3] = 1
$data[0,4] = 'LP1912'
$data[1,0] = '13:33:42'
$data[1,1] = '13:34'
$data[1,2] = '23_HERNANDEZ'
$data[1,3] = 1
$data[1,4] = 'LP1912'
$ws1.Range("A233:E234").Value = $data

$data = New-Object 'object[,]' 2,5
$data[0,0] = '12:59:47'
$data[0,1] = '13:56'
$data[0,2] = '16_P MOR-167 Y 521'
$data[0,3] = 57
$data[0,4] = 'LP1912'
$data[1,0] = '12:21:08'
$data[1,1] = '13:56'
$data[1,2] = '225_GOMEZ'
$data[1,3] = 95
$data[1,4] = 'LP1912'
$ws1.Range("A242:E243").Value = $data

$data = New-Object 'object[,]' 34,5
$data[0,0] = '14:45:17'
$data[0,1] = '14:45'
$data[0,2] = '15_ABASTO'
$data[0,3] = 0
$data[0,4] = 'LP1912'
$data[1,0] = '12:47:27'
$data[1,1] = '14:45'
$data[1,2] = '14_ABASTO'
$data[1,3] = 118
$data[1,4] = 'LP1912'
$data[2,0] = '14:45:17'
$data[2,1] = '14:46'
$data[2,2] = '16_SANTA ANA'
$data[2,3] = 1
$data[2,4] = 'LP1912'
$data[3,0] = '12:59:47'
$data[3,1] = '14:56'
$data[3,2] = '16_P MOR-SANTA ANA'
$data[3,3] = 117
$data[3,4] = 'LP1912'
$data[4,0] = '13:59:06'
$data[4,1] = '14:57'
$data[4,2] = '16_P MOR-SANTA ANA'
$data[4,3] = 58
$data[4,4] = 'LP1912'
$data[5,0] = '12:59:47'
$data[5,1] = '14:58'
$data[5,2] = '215B_EL PATO'
$data[5,3] = 119
$data[5,4] = 'LP1912'
$data[6,0] = '13:33:42'
$data[6,1] = '15:00'
$data[6,2] = '81_EL PELIGRO'
$data[6,3] = 87
$data[6,4] = 'LP1912'
$data[7,0] = '13:33:42'
$data[7,1] = '15:05'
$data[7,2] = '10_OLMOS'
$data[7,3] = 92
$data[7,4] = 'LP1912'
$data[8,0] = '14:45:17'
$data[8,1] = '15:05'
$data[8,2] = '23_HERNANDEZ'
$data[8,3] = 20
$data[8,4] = 'LP1912'
$data[9,0] = '13:59:06'
$data[9,1] = '15:10'
$data[9,2] = '17_ROMERO'
$data[9,3] = 71
$data[9,4] = 'LP1912'
$data[10,0] = '13:33:42'
$data[10,1] = '15:13'
$data[10,2] = '11_ETCHEVERRY'
$data[10,3] = 100
$data[10,4] = 'LP1912'
$data[11,0] = '13:59:06'
$data[11,1] = '15:14'
$data[11,2] = '11_ETCHEVERRY'
$data[11,3] = 75
$data[11,4] = 'LP1912'
$data[12,0] = '13:33:42'
$data[12,1] = '15:17'
$data[12,2] = '26_HERNANDEZ'
$data[12,3] = 104
$data[12,4] = 'LP1912'
$data[13,0] = '13:59:06'
$data[13,1] = '15:18'
$data[13,2] = '26_HERNANDEZ'
$data[13,3] = 79
$data[13,4] = 'LP1912'
$data[14,0] = '14:24:16'
$data[14,1] = '15:21'
$data[14,2] = '26_HERNANDEZ'
$data[14,3] = 57
$data[14,4] = 'LP1912'
$data[15,0] = '14:24:16'
$data[15,1] = '15:32'
$data[15,2] = '84_COLONIA URQUIZA-ESC 49'
$data[15,3] = 68
$data[15,4] = 'LP1912'
$data[16,0] = '13:59:06'
$data[16,1] = '15:35'
$data[16,2] = '23_HERNANDEZ'
$data[16,3] = 96
$data[16,4] = 'LP1912'
$data[17,0] = '13:59:06'
$data[17,1] = '15:37'
$data[17,2] = '10_OLMOS'
$data[17,3] = 98
$data[17,4] = 'LP1912'
$data[18,0] = '14:24:16'
$data[18,1] = '15:38'
$data[18,2] = '23_HERNANDEZ'
$data[18,3] = 74
$data[18,4] = 'LP1912'
$data[19,0] = '14:45:17'
$data[19,1] = '15:38'
$data[19,2] = '215A_EL PATO'
$data[19,3] = 53
$data[19,4] = 'LP1912'
$data[20,0] = '13:59:06'
$data[20,1] = '15:39'
$data[20,2] = '215A_EL PATO'
$data[20,3] = 100
$data[20,4] = 'LP1912'
$data[21,0] = '14:24:16'
$data[21,1] = '15:46'
$data[21,2] = '14_ABASTO'
$data[21,3] = 82
$data[21,4] = 'LP1912'
$data[22,0] = '13:59:06'
$data[22,1] = '15:47'
$data[22,2] = '16_P MOR-167 Y 521'
$data[22,3] = 108
$data[22,4] = 'LP1912'
$data[23,0] = '13:59:06'
$data[23,1] = '15:48'
$data[23,2] = '14_ABASTO'
$data[23,3] = 109
$data[23,4] = 'LP1912'
$data[24,0] = '13:59:06'
$data[24,1] = '15:54'
$data[24,2] = '11_ETCHEVERRY'
$data[24,3] = 115
$data[24,4] = 'LP1912'
$data[25,0] = '14:24:16'
$data[25,1] = '15:56'
$data[25,2] = '17_ROMERO'
$data[25,3] = 92
$data[25,4] = 'LP1912'
$data[26,0] = '13:59:06'
$data[26,1] = '15:57'
$data[26,2] = '27_EL RETIRO'
$data[26,3] = 118
$data[26,4] = 'LP1912'
$data[27,0] = '14:45:17'
$data[27,1] = '16:09'
$data[27,2] = '14_ABASTO'
$data[27,3] = 84
$data[27,4] = 'LP1912'
$data[28,0] = '14:24:16'
$data[28,1] = '16:15'
$data[28,2] = '225_C ROCA-H SUR'
$data[28,3] = 111
$data[28,4] = 'LP1912'
$data[29,0] = '14:24:16'
$data[29,1] = '16:20'
$data[29,2] = '215C_EL PATO'
$data[29,3] = 116
$data[29,4] = 'LP1912'
$data[30,0] = '14:24:16'
$data[30,1] = '16:21'
$data[30,2] = '26_HERNANDEZ'
$data[30,3] = 117
$data[30,4] = 'LP1912'
$data[31,0] = '14:45:17'
$data[31,1] = '16:30'
$data[31,2] = '15_ABASTO'
$data[31,3] = 105
$data[31,4] = 'LP1912'
$data[32,0] = '14:45:17'
$data[32,1] = '16:43'
$data[32,2] = '16_P MOR-SANTA ANA'
$data[32,3] = 118
$data[32,4] = 'LP1912'
$data[33,0] = '14:45:17'
$data[33,1] = '16:43'
$data[33,2] = '225_GOMEZ'
$data[33,3] = 118
$data[33,4] = 'LP1912'
$ws1.Range("A264:E297").Value = $data

# ---- Sheet: LP1912-215 ----
$data = New-Object 'object[,]' 2,5
$data[0,0] = 'Última actualización: 14:45:17'
$data[1,0] = 'Total filas: 33'
$ws2.Range("A2:E3").Value = $data

$data = New-Object 'object[,]' 3,5
$data[0,0] = '14:45:17'
$data[0,1] = '15:38'
$data[0,2] = '215A_EL PATO'
$data[0,3] = 53
$data[0,4] = 'LP1912'
$data[1,0] = '13:59:06'
$data[1,1] = '15:39'
$data[1,2] = '215A_EL PATO'
$data[1,3] = 100
$data[1,4] = 'LP1912'
$data[2,0] = '14:24:16'
$data[2,1] = '16:20'
$data[2,2] = '215C_EL PATO'
$data[2,3] = 116
$ws2.Range("A36:E38").Value = $data

# ---- Sheet: 6203-6173 ----
$data = New-Object 'object[,]' 1,5
$data[0,0] = 'Última actualización: 14:45:17'
$ws3.Range("A2:E2").Value = $data
